$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("B2").Value = 0.2122186495176849
$ws.Range("C2").Value = 0.5112540192926045
$ws.Range("J2").Value = 0.01286173633440514
$ws.Range("P2").Value = 0.1607717041800643
$ws.Range("S2").Value = 0.1028938906752412

# Row 3
$ws.Range("B3").Value = 0.006060606060606061
$ws.Range("C3").Value = 0.04242424242424243
$ws.Range("J3").Value = 0.0303030303030303
$ws.Range("P3").Value = 0.6787878787878788
$ws.Range("S3").Value = 0.2424242424242424

# Row 4
$ws.Range("J4").Value = 0.1052631578947368
$ws.Range("P4").Value = 0.7105263157894737
$ws.Range("S4").Value = 0.1842105263157895

# Row 6
$ws.Range("B6").Value = 0.04824561403508772
$ws.Range("D6").Value = 0.02192982456140351
$ws.Range("F6").Value = 0.04824561403508772
$ws.Range("J6").Value = 0.2675438596491228
$ws.Range("O6").Value = 0.008771929824561403
$ws.Range("Q6").Value = 0.1666666666666667
$ws.Range("R6").Value = 0.08333333333333333
$ws.Range("S6").Value = 0.3552631578947368

# Row 7
$ws.Range("B7").Value = 0.1479289940828402
$ws.Range("D7").Value = 0.005917159763313609
$ws.Range("F7").Value = 0.05917159763313609
$ws.Range("J7").Value = 0.1479289940828402
$ws.Range("Q7").Value = 0.1301775147928994
$ws.Range("R7").Value = 0.05917159763313609
$ws.Range("S7").Value = 0.4497041420118343

# Row 8
$ws.Range("B8").Value = 0.1005154639175258
$ws.Range("D8").Value = 0.01804123711340206
$ws.Range("E8").Value = 0.002577319587628866
$ws.Range("F8").Value = 0.06185567010309279
$ws.Range("J8").Value = 0.1030927835051546
$ws.Range("O8").Value = 0.01288659793814433
$ws.Range("Q8").Value = 0.1675257731958763
$ws.Range("R8").Value = 0.1056701030927835
$ws.Range("S8").Value = 0.4278350515463917

# Row 9
$ws.Range("B9").Value = 0.07960199004975124
$ws.Range("D9").Value = 0.03482587064676617
$ws.Range("F9").Value = 0.07960199004975124
$ws.Range("J9").Value = 0.05970149253731343
$ws.Range("O9").Value = 0.03482587064676617
$ws.Range("Q9").Value = 0.2537313432835821
$ws.Range("R9").Value = 0.06467661691542288
$ws.Range("S9").Value = 0.3930348258706468

# Row 10
$ws.Range("B10").Value = 0.1213114754098361
$ws.Range("D10").Value = 0.01557377049180328
$ws.Range("E10").Value = 0.000819672131147541
$ws.Range("F10").Value = 0.07459016393442623
$ws.Range("J10").Value = 0.08442622950819673
$ws.Range("O10").Value = 0.01967213114754099
$ws.Range("Q10").Value = 0.1704918032786885
$ws.Range("R10").Value = 0.0959016393442623
$ws.Range("S10").Value = 0.4172131147540983

# Row 11
$ws.Range("F11").Value = 0.003546099290780142
$ws.Range("G11").Value = 0.1347517730496454
$ws.Range("J11").Value = 0.124113475177305
$ws.Range("K11").Value = 0.198581560283688
$ws.Range("L11").Value = 0.5212765957446809
$ws.Range("S11").Value = 0.01773049645390071

# Row 12
$ws.Range("G12").Value = 0.7837837837837838
$ws.Range("J12").Value = 0.1621621621621622
$ws.Range("K12").Value = 0.01351351351351351
$ws.Range("S12").Value = 0.04054054054054054

# Row 13
$ws.Range("G13").Value = 0.5517241379310345
$ws.Range("J13").Value = 0.3103448275862069
$ws.Range("S13").Value = 0.1379310344827586

# Row 15
$ws.Range("F15").Value = 0.05092592592592592
$ws.Range("H15").Value = 0.1527777777777778
$ws.Range("I15").Value = 0.1064814814814815
$ws.Range("J15").Value = 0.3564814814814815
$ws.Range("K15").Value = 0.05555555555555555
$ws.Range("M15").Value = 0.01388888888888889
$ws.Range("O15").Value = 0.02777777777777778
$ws.Range("S15").Value = 0.2361111111111111

# Row 16
$ws.Range("F16").Value = 0.02747252747252747
$ws.Range("H16").Value = 0.1703296703296703
$ws.Range("I16").Value = 0.09340659340659341
$ws.Range("J16").Value = 0.3791208791208791
$ws.Range("K16").Value = 0.1208791208791209
$ws.Range("M16").Value = 0.02197802197802198
$ws.Range("O16").Value = 0.06593406593406594
$ws.Range("S16").Value = 0.1208791208791209

# Row 17
$ws.Range("F17").Value = 0.005277044854881266
$ws.Range("H17").Value = 0.1609498680738786
$ws.Range("I17").Value = 0.09762532981530343
$ws.Range("J17").Value = 0.4722955145118733
$ws.Range("K17").Value = 0.079155672823219
$ws.Range("M17").Value = 0.0158311345646438
$ws.Range("N17").Value = 0.002638522427440633
$ws.Range("O17").Value = 0.05277044854881267
$ws.Range("S17").Value = 0.1134564643799472

# Row 18
$ws.Range("F18").Value = 0.0202020202020202
$ws.Range("H18").Value = 0.1565656565656566
$ws.Range("I18").Value = 0.0707070707070707
$ws.Range("J18").Value = 0.4747474747474748
$ws.Range("K18").Value = 0.0707070707070707
$ws.Range("M18").Value = 0.005050505050505051
$ws.Range("O18").Value = 0.0707070707070707
$ws.Range("S18").Value = 0.1313131313131313

# Row 19
$ws.Range("F19").Value = 0.02199528672427337
$ws.Range("H19").Value = 0.1846032992930086
$ws.Range("I19").Value = 0.08876669285153181
$ws.Range("J19").Value = 0.3935585231736057
$ws.Range("K19").Value = 0.1076197957580518
$ws.Range("M19").Value = 0.01256873527101336
$ws.Range("N19").Value = 0.001571091908876669
$ws.Range("O19").Value = 0.07855459544383346
$ws.Range("S19").Value = 0.1107619795758052
